$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Q1" (column D) forecast-error values for rows 113-116
$ws.Range("D113").Value = 0.7254492243564907
$ws.Range("D114").Value = 0.7215746373564907
$ws.Range("D115").Value = 0.5311946523564907
$ws.Range("D116").Value = 0.5539812373564907

# Add the "Q0" (column C) forecast-error values for rows 117-122
$ws.Range("C117").Value = 0.1753415943564907
$ws.Range("C118").Value = 0.2651053283564908
$ws.Range("C119").Value = 0.08763596535649075
$ws.Range("C120").Value = 0.1003532183564907
$ws.Range("C121").Value = -0.02418658464350926
$ws.Range("C122").Value = 0.2001520573564908
